# "Refined metadata to be additional tab"
#
# 1) Update the "panel_query_time" (F column) timestamps on the existing
#    "data" sheet to reflect a fresh query run.
# 2) Add a new "metadata" worksheet (after "data") carrying one summary
#    row describing the panel query itself (name/id/version/etc.).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "data"

# --- 1) refresh the per-gene query timestamps in column F -------------
$ws1.Range("F2").Value = "2021-10-05 14:33:44.188007"
$ws1.Range("F3").Value = "2021-10-05 14:33:44.188015"
$ws1.Range("F4").Value = "2021-10-05 14:33:44.188018"
$ws1.Range("F5").Value = "2021-10-05 14:33:44.188020"
$ws1.Range("F6").Value = "2021-10-05 14:33:44.188023"
$ws1.Range("F7").Value = "2021-10-05 14:33:44.188026"
$ws1.Range("F8").Value = "2021-10-05 14:33:44.188028"
$ws1.Range("F9").Value = "2021-10-05 14:33:44.188031"
$ws1.Range("F10").Value = "2021-10-05 14:33:44.188034"
$ws1.Range("F11").Value = "2021-10-05 14:33:44.188036"
$ws1.Range("F12").Value = "2021-10-05 14:33:44.188039"
$ws1.Range("F13").Value = "2021-10-05 14:33:44.188041"
$ws1.Range("F14").Value = "2021-10-05 14:33:44.188044"
$ws1.Range("F15").Value = "2021-10-05 14:33:44.188046"
$ws1.Range("F16").Value = "2021-10-05 14:33:44.188049"
$ws1.Range("F17").Value = "2021-10-05 14:33:44.188052"
$ws1.Range("F18").Value = "2021-10-05 14:33:44.188054"
$ws1.Range("F19").Value = "2021-10-05 14:33:44.188057"
$ws1.Range("F20").Value = "2021-10-05 14:33:44.188060"
$ws1.Range("F21").Value = "2021-10-05 14:33:44.188062"
$ws1.Range("F22").Value = "2021-10-05 14:33:44.188065"
$ws1.Range("F23").Value = "2021-10-05 14:33:44.188067"
$ws1.Range("F24").Value = "2021-10-05 14:33:44.188070"
$ws1.Range("F25").Value = "2021-10-05 14:33:44.188072"
$ws1.Range("F26").Value = "2021-10-05 14:33:44.188075"
$ws1.Range("F27").Value = "2021-10-05 14:33:44.188078"
$ws1.Range("F28").Value = "2021-10-05 14:33:44.188080"
$ws1.Range("F29").Value = "2021-10-05 14:33:44.188082"
$ws1.Range("F30").Value = "2021-10-05 14:33:44.188085"
$ws1.Range("F31").Value = "2021-10-05 14:33:44.188087"
$ws1.Range("F32").Value = "2021-10-05 14:33:44.188090"
$ws1.Range("F33").Value = "2021-10-05 14:33:44.188092"
$ws1.Range("F34").Value = "2021-10-05 14:33:44.188095"
$ws1.Range("F35").Value = "2021-10-05 14:33:44.188098"
$ws1.Range("F36").Value = "2021-10-05 14:33:44.188100"
$ws1.Range("F37").Value = "2021-10-05 14:33:44.188103"
$ws1.Range("F38").Value = "2021-10-05 14:33:44.188105"
$ws1.Range("F39").Value = "2021-10-05 14:33:44.188108"
$ws1.Range("F40").Value = "2021-10-05 14:33:44.188111"
$ws1.Range("F41").Value = "2021-10-05 14:33:44.188113"
$ws1.Range("F42").Value = "2021-10-05 14:33:44.188116"
$ws1.Range("F43").Value = "2021-10-05 14:33:44.188119"
$ws1.Range("F44").Value = "2021-10-05 14:33:44.188121"
$ws1.Range("F45").Value = "2021-10-05 14:33:44.188124"
$ws1.Range("F46").Value = "2021-10-05 14:33:44.188126"
$ws1.Range("F47").Value = "2021-10-05 14:33:44.188129"
$ws1.Range("F48").Value = "2021-10-05 14:33:44.188132"
$ws1.Range("F49").Value = "2021-10-05 14:33:44.188134"
$ws1.Range("F50").Value = "2021-10-05 14:33:44.188137"
$ws1.Range("F51").Value = "2021-10-05 14:33:44.188139"
$ws1.Range("F52").Value = "2021-10-05 14:33:44.188142"

# --- 2) add the "metadata" worksheet, placed after "data" -------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# header row styling (bold / bordered / centered) copied from "data"!B1:F1
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Dilated Cardiomyopathy"
$ws2.Range("C2").Value = 95

$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "1.4"
$ws2.Range("D2").Style = "Normal"

$ws2.Range("E2").Value = "2021-08-19T21:18:29.514177Z"
$ws2.Range("F2").Value = "2021-10-05 14:33:44.184366"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/95/?format=json"

[void]$ws1.Select()
[void]$ws1.Range("A1").Select()
